$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  D = 10;                  E = 9.710000038146973;  F = 10;                  G = 9.560000419616699;  H = 436518932 },
    @{ Row = 3;  D = 10;                  E = 9.710000038146973;  F = 10;                  G = 9.560000419616699;  H = 436518932 },
    @{ Row = 4;  D = 9.850000381469728;   E = 9.779999732971191;  F = 9.85200023651123;    G = 9.760000228881836;  H = 436518932 },
    @{ Row = 5;  D = 9.890000343322754;   E = 9.970000267028809;  F = 9.970000267028809;   G = 9.829999923706056;  H = 436518932 },
    @{ Row = 6;  D = 5.659999847412109;   E = 5.110000133514404;  F = 6.091000080108643;   G = 4.449999809265137;  H = 436518932 },
    @{ Row = 7;  D = 8.819999694824219;   E = 10.01000022888184;  F = 16.1200008392334;    G = 8.039999961853027;  H = 436518932 },
    @{ Row = 8;  D = 3.670000076293945;   E = 2.539999961853028;  F = 3.839999914169312;   G = 2.440000057220459;  H = 436518932 },
    @{ Row = 9;  D = 1.590000033378601;   E = 1.350000023841858;  F = 2.079999923706055;   G = 1.110000014305115;  H = 436518932 },
    @{ Row = 10; D = 0.7160000205039978;  E = 3.269999980926514;  F = 4.449999809265137;   G = 0.6740000247955322; H = 436518932 },
    @{ Row = 11; D = 2.5;                 E = 2.900000095367432;  F = 3.740000009536743;   G = 2.019999980926514;  H = 436518932 },
    @{ Row = 12; D = 2.394999980926514;   E = 2.009999990463257;  F = 2.430000066757202;   G = 1.740000009536743;  H = 436518932 },
    @{ Row = 13; D = 1.549999952316284;   E = 1.269999980926514;  F = 1.590000033378601;   G = 1.159999966621399;  H = 436518932 },
    @{ Row = 14; D = 2.069999933242798;   E = 1.620000004768372;  F = 2.240000009536743;   G = 1.620000004768372;  H = 436518932 },
    @{ Row = 15; D = 2.039999961853028;   E = 1.659999966621399;  F = 2.069999933242798;   G = 1.440000057220459;  H = 436518932 },
    @{ Row = 16; D = 1.490000009536743;   E = 1.509999990463257;  F = 1.669999957084656;   G = 1.370000004768372;  H = 436518932 },
    @{ Row = 17; D = 1.440000057220459;   E = 1.590000033378601;  F = 1.940000057220459;   G = 1.360000014305115;  H = 436518932 },
    @{ Row = 18; D = 4.585000038146973;   E = 4.239999771118164;  F = 4.820000171661377;   G = 2.849999904632568;  H = 436518932 },
    @{ Row = 19; D = 2.900000095367432;   E = 3.410000085830688;  F = 3.789999961853027;   G = 2.359999895095825;  H = 436518932 },
    @{ Row = 20; D = 7.150000095367432;   E = 6.349999904632568;  F = 8.720000267028809;   G = 6.090000152587891;  H = 436518932 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = "BBAI"
}
